# Daily attendance processing - 2026-01-12 13:56:58
# Reorders the "Recorded By" (column G) values for specific known
# author combinations so that the system account is listed after the
# human/automated account, instead of before it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Map of exact current cell text -> replacement cell text.
$replacements = @{
    "System, system, backup@backdoor.com" = "System, backup@backdoor.com, system"
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "System, admin@admin.com"             = "admin@admin.com, System"
    "admin@admin.com, dnasr281@gmail.com" = "dnasr281@gmail.com, admin@admin.com"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value()
    if ($null -ne $val -and $replacements.ContainsKey($val)) {
        $cell.Value = $replacements[$val]
    }
}
